$d = $word.ActiveDocument

# The last paragraph in the document currently holds only the "晴转多云" run
# (it also carries a <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
# paragraph mark). We need to:
#   1) Insert four new runs ("0", "6:", "20-08:30", "晴空万里") in front of
#      that run, inside the SAME paragraph (so it keeps its original
#      w14:paraId / rsid attributes and its <w:pPr>).
#   2) Push the existing "晴转多云" run out into a brand new, attribute-less
#      <w:p> paragraph that follows.
#
# Range.InsertXML lets us replace a Range's contents with arbitrary OOXML,
# including full paragraph marks, in one shot - which is the only way to
# get distinct <w:r> elements (some hinted eastAsia, some not) without Word
# silently coalescing them back into a single run.

$target = $d.Paragraphs.Item($d.Paragraphs.Count)
$pid = $target.Range.Start
$startPara = $target.Range.End

$r = $d.Range($pid, $startPara)

$xml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
          '<w:body>' +
            '<w:p w14:paraId="7D73B72B" w14:textId="2E386EA2" w:rsidR="005E23FA" w:rsidRDefault="005E23FA">' +
              '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' +
              '<w:r><w:t>0</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>6:</w:t></w:r>' +
              '<w:r><w:t>20-08:30</w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>晴空万里</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>晴转多云</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($xml)
